# CIV-14009 Template and variables update
#
# 1) Body text: placeholder "paidInFullDate" -> "dateOfEvent"
#    (keeps the surrounding "<<" / ">>" runs and the spellStart/spellEnd
#     proofErr markers untouched).
# 2) Footer: "((Name and Contact details of the owning court))" is
#    replaced by the merge-field style placeholder
#    "<<nameAndContactDetailsOwningCourt>>", typed out as a sequence of
#    small runs (mirroring how Word splits camelCase words while the
#    spell-checker runs), wrapped in a proofErr spellStart/spellEnd pair.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: paidInFullDate -> dateOfEvent
# ---------------------------------------------------------------------
# Find the run together with everything that follows it up to (and
# including) the end of the paragraph, so the in-place XML splice lands
# correctly instead of being appended after the paragraph's last run.
$body = $d.Content
$found1 = $body.Find.Execute("paidInFullDate>>. ")
if ($found1) {
    $tailRange = $d.Range($body.Start, $body.End)

    $xml1 = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r w:rsidR="006D5997"><w:t>dateOfEvent</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r w:rsidR="00F751E5"><w:t>&gt;&gt;</w:t></w:r>' +
        '<w:r w:rsidR="00BF745E" w:rsidRPr="00BF745E"><w:t xml:space="preserve">. </w:t></w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $tailRange.InsertXML($xml1)
}

# ---------------------------------------------------------------------
# Edit 2: ((Name and Contact details of the owning court))
#         -> <<nameAndContactDetailsOwningCourt>>
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftrRange = $ftr.Range
$found2 = $ftrRange.Find.Execute("((Name and Contact details of the owning court))")

if ($found2) {
    $xml2 = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r><w:t>&lt;&lt;</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>n</w:t></w:r>' +
        '<w:r><w:t>ame</w:t></w:r>' +
        '<w:r><w:t>AndC</w:t></w:r>' +
        '<w:r><w:t>ontact</w:t></w:r>' +
        '<w:r><w:t>D</w:t></w:r>' +
        '<w:r><w:t>etails</w:t></w:r>' +
        '<w:r><w:t>O</w:t></w:r>' +
        '<w:r><w:t>wning</w:t></w:r>' +
        '<w:r><w:t>C</w:t></w:r>' +
        '<w:r><w:t>ourt</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>&gt;&gt;</w:t></w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $ftrRange.InsertXML($xml2)
}
